$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode values:
# Row 2 (A Suite / IAM module): Y -> N
# Row 4 (C Suite / Authoring module): N -> Y
# Row 5 (D Suite / Profile module): N -> Y
$ws.Range("C2").Value = "N"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"

# Update the active cell selection to C5
$ws.Range("C5").Select()
